$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.878.34"
$ws.Range("E2").Value = "  -1.58%  "

# Row 3
$ws.Range("D3").Value = "3.383.18"
$ws.Range("E3").Value = "  -0.90%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.27"
$ws.Range("E5").Value = "  -0.90%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.60"
$ws.Range("E6").Value = "  -1.51%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("D8").Value = "3.381.25"
$ws.Range("E8").Value = "  -0.95%  "

# Row 9
$ws.Range("E9").Value = "  -1.42%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.53"
$ws.Range("E10").Value = "  +0.52%  "

# Row 11
$ws.Range("E11").Value = "  -3.31%  "

# Row 12
$ws.Range("E12").Value = "  -2.08%  "

# Row 13
$ws.Range("D13").Value = "3.957.07"
$ws.Range("E13").Value = "  -1.04%  "

# Row 14
$ws.Range("E14").Value = "  +1.10%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.37"
$ws.Range("E15").Value = "  +3.52%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000172"
$ws.Range("E16").Value = "  -4.63%  "

# Row 17
$ws.Range("D17").Value = "3.379.94"
$ws.Range("E17").Value = "  -0.92%  "

# Row 18
$ws.Range("D18").Value = "61.064.65"
$ws.Range("E18").Value = "  -1.43%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.00"
$ws.Range("E19").Value = "  -0.87%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.83"
$ws.Range("E20").Value = "  -1.29%  "

# Row 21
$ws.Range("E21").Value = "  -1.45%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "376.80"
$ws.Range("E22").Value = "  -3.63%  "

# Row 23
$ws.Range("E23").Value = "  -2.84%  "

# Row 24
$ws.Range("D24").Value = "3.508.28"
$ws.Range("E24").Value = "  -1.22%  "

# Row 25
$ws.Range("E25").Value = "  +0.11%  "

# Row 26
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.24"
$ws.Range("E26").Value = "  -0.60%  "

# Row 27
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000125"
$ws.Range("E27").Value = "  -2.64%  "

# Row 28
$ws.Range("E28").Value = "  +11.35%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.53"
$ws.Range("E29").Value = "  -1.89%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.168"
$ws.Range("E30").Value = "  +4.60%  "

# Row 31
$ws.Range("E31").Value = "  -0.10%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.15"
$ws.Range("E32").Value = "  -1.99%  "

# Row 33
$ws.Range("E33").Value = "  -0.91%  "

# Row 34
$ws.Range("E34").Value = "  -0.01%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.76"
$ws.Range("E35").Value = "  +0.76%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.21"
$ws.Range("E36").Value = "  -5.15%  "

# Row 37
$ws.Range("E37").Value = "  -2.67%  "

# Row 38
$ws.Range("E38").Value = "  -2.09%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "164.75"
$ws.Range("E39").Value = "  +1.24%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0755"
$ws.Range("E40").Value = "  -4.69%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.06%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.771"
$ws.Range("E42").Value = "  -2.70%  "

# Row 43
$ws.Range("E43").Value = "  -3.77%  "

# Row 44
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.54"
$ws.Range("E44").Value = "  -0.27%  "

# Row 45
$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.21"
$ws.Range("E45").Value = "  -1.88%  "

# Row 46
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.40"
$ws.Range("E46").Value = "  -1.70%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.03"
$ws.Range("E47").Value = "  -5.08%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.39"
$ws.Range("E48").Value = "  +0.82%  "

# Row 49
$ws.Range("E49").Value = "  -2.62%  "

# Row 50
$ws.Range("D50").Value = "2.379.92"
$ws.Range("E50").Value = "  -0.56%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.40"
$ws.Range("E51").Value = "  +4.21%  "
